$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (A1:G18), values kept as text since original cells are inlineStr.
$data = @(
    @("core_id","er.c1","er.c2","er.c3","p53.c1","p53.c2","p53.c3"),
    @("1","9","9","9","3","0","5"),
    @("4","2","0","0","8","4","4"),
    @("10","x","9","0","3","5","5"),
    @("13","9","1","1","5","1","3"),
    @("3","2","x","1","5","9","9"),
    @("12","9","x","x","1","2","4"),
    @("5","2","2","2","3","2","1"),
    @("14","1","2","1","3","4","5"),
    @("2","1","2","2","8","9","x"),
    @("11","2","x","0","8","9","5"),
    @("6","x","2","2","5","3","3"),
    @("9","0","1","x","3","3","2"),
    @("15","1","1","2","3","x","1"),
    @("8","0","1","x","x","9","3"),
    @("17","0","1","2","8","3","2"),
    @("7","x","1","x","4","4","5"),
    @("16","x","2","","8","4","")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        if ($value -eq "") {
            $cell.Value = $null
        } else {
            $cell.Value = $value
        }
    }
}

# Make sure the new header cells (F1:G1) match the bold/centered style used by
# the rest of the header row (A1:E1).
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").HorizontalAlignment = -4108  # xlCenter
